# Update countries & provincias Spain
# - Re-sort three countries (Republica de Macedonia, Kenia, Togo) to new
#   positions in the list and refresh their statistics.
# - Refresh several other countries' statistics in place.
# - Bump the "datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow {
    param(
        [int]$Row,
        [int]$B, [int]$C, [int]$D, [int]$E, [int]$F, [int]$G, [int]$H
    )
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
}

# ---------------------------------------------------------------------
# 1. Simple in-place statistic refreshes (no row movement)
# ---------------------------------------------------------------------

# Paises Bajos (row 15)
$ws.Cells.Item(15, 6).Value = 1182

# Israel (row 21)
$ws.Cells.Item(21, 5).Value = 6653
$ws.Cells.Item(21, 7).Value = 3
$ws.Cells.Item(21, 8).Value = 39

# Australia (row 23)
$ws.Cells.Item(23, 6).Value = 85

# Finlandia (row 42)
$ws.Cells.Item(42, 5).Value = 1294
$ws.Cells.Item(42, 6).Value = 68
$ws.Cells.Item(42, 7).Value = 2
$ws.Cells.Item(42, 8).Value = 21

# Bosnia y Herzegovina (row 71)
$ws.Cells.Item(71, 2).Value = 561
$ws.Cells.Item(71, 3).Value = 28
$ws.Cells.Item(71, 5).Value = 517
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 17

# Sri Lanka (row 109)
$ws.Cells.Item(109, 2).Value = 152
$ws.Cells.Item(109, 3).Value = 1
$ws.Cells.Item(109, 4).Value = 24
$ws.Cells.Item(109, 5).Value = 124

# ---------------------------------------------------------------------
# 2. Republica de Macedonia: move from its old spot (just above Uruguay)
#    to a new spot right before Kuwait, with refreshed figures.
# ---------------------------------------------------------------------
$oldRow = 0
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    if ($ws.Cells.Item($r, 1).Text -eq "Republica de Macedonia") {
        $oldRow = $r
        break
    }
}
$ws.Rows($oldRow).Delete()

$newRow = 0
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    if ($ws.Cells.Item($r, 1).Text -eq "Kuwait") {
        $newRow = $r
        break
    }
}
$ws.Rows($newRow).Insert()
$ws.Cells.Item($newRow, 1).Value = "Republica de Macedonia"
Set-CountryRow $newRow 430 46 17 402 8 0 11

# ---------------------------------------------------------------------
# 3. Kenia: move from its old spot (just above Niger) to a new spot
#    right before Mayotte, with refreshed figures.
# ---------------------------------------------------------------------
$oldRow = 0
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    if ($ws.Cells.Item($r, 1).Text -eq "Kenia") {
        $oldRow = $r
        break
    }
}
$ws.Rows($oldRow).Delete()

$newRow = 0
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    if ($ws.Cells.Item($r, 1).Text -eq "Mayotte") {
        $newRow = $r
        break
    }
}
$ws.Rows($newRow).Insert()
$ws.Cells.Item($newRow, 1).Value = "Kenia"
Set-CountryRow $newRow 122 12 4 114 2 1 4

# ---------------------------------------------------------------------
# 4. Togo: move from its old spot (just above Polinesia Francesa) to a
#    new spot right before Puerto Rico, with refreshed figures.
# ---------------------------------------------------------------------
$oldRow = 0
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    if ($ws.Cells.Item($r, 1).Text -eq "Togo") {
        $oldRow = $r
        break
    }
}
$ws.Rows($oldRow).Delete()

$newRow = 0
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    if ($ws.Cells.Item($r, 1).Text -eq "Puerto Rico") {
        $newRow = $r
        break
    }
}
$ws.Rows($newRow).Insert()
$ws.Cells.Item($newRow, 1).Value = "Togo"
Set-CountryRow $newRow 40 1 17 20 0 1 3

# ---------------------------------------------------------------------
# 5. Bump the "datos actualizados" timestamp in A1.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 14:50"
